# Add HADDOCK results for BA.2.75 and BA.5 for CC12.1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (BA.2.75) - HADDOCK metrics, columns K..Z
$ws.Range("K2").Value = -118.8
$ws.Range("L2").Value = 10.5
$ws.Range("M2").Value = 26
$ws.Range("N2").Value = 0.9
$ws.Range("O2").Value = 0.6
$ws.Range("P2").Value = -72.3
$ws.Range("Q2").Value = 8
$ws.Range("R2").Value = -265.7
$ws.Range("S2").Value = 48.8
$ws.Range("T2").Value = -8.9
$ws.Range("U2").Value = 2.3
$ws.Range("V2").Value = 155.7
$ws.Range("W2").Value = 53.5
$ws.Range("X2").Value = 2265.6
$ws.Range("Y2").Value = 134.3
$ws.Range("Z2").Value = -1.6

# Row 3 (BA.5) - HADDOCK metrics, columns K..Z
$ws.Range("K3").Value = -121.4
$ws.Range("L3").Value = 3.2
$ws.Range("M3").Value = 18
$ws.Range("N3").Value = 0.6
$ws.Range("O3").Value = 0.4
$ws.Range("P3").Value = -78
$ws.Range("Q3").Value = 3.1
$ws.Range("R3").Value = -198.3
$ws.Range("S3").Value = 16.9
$ws.Range("T3").Value = -18.6
$ws.Range("U3").Value = 2.9
$ws.Range("V3").Value = 147.3
$ws.Range("W3").Value = 36.6
$ws.Range("X3").Value = 2326.3
$ws.Range("Y3").Value = 67.6
$ws.Range("Z3").Value = -2.2

# Selection moved to L7 as last-saved cursor position
$ws.Range("L7").Select() | Out-Null
